$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($Ws, $Row, $Col, $Text) {
    $cell = $Ws.Cells.Item($Row, $Col)
    # Force the cell to Text format first so Excel does not reinterpret
    # numeric-looking strings (e.g. "1.003", "0.05380", "27.922.19") as
    # numbers/dates and strip meaningful digits (like trailing zeros).
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

Set-CellText $ws 2 4 "27.922.19"
Set-CellText $ws 2 5 "  -4.35%  "
Set-CellText $ws 3 4 "1.739.74"
Set-CellText $ws 3 5 "  -4.71%  "
Set-CellText $ws 5 4 "226.64"
Set-CellText $ws 5 5 "  -3.53%  "
Set-CellText $ws 6 4 "0.5785"
Set-CellText $ws 6 5 "  -3.51%  "
Set-CellText $ws 7 4 "1.003"
Set-CellText $ws 8 4 "0.2741"
Set-CellText $ws 8 5 "  -0.87%  "
Set-CellText $ws 9 4 "23.11"
Set-CellText $ws 9 5 "  -1.49%  "
Set-CellText $ws 10 4 "0.06632"
Set-CellText $ws 10 5 "  -4.37%  "
Set-CellText $ws 11 4 "0.07558"
Set-CellText $ws 11 5 "  -0.67%  "
Set-CellText $ws 12 4 "1.744.41"
Set-CellText $ws 12 5 "  -4.40%  "
Set-CellText $ws 13 5 "  -0.28%  "
Set-CellText $ws 14 4 "0.6026"
Set-CellText $ws 14 5 "  -3.75%  "
Set-CellText $ws 15 4 "1.976.98"
Set-CellText $ws 15 5 "  -4.65%  "
Set-CellText $ws 16 4 "74.57"
Set-CellText $ws 16 5 "  -3.62%  "
Set-CellText $ws 17 4 "0.000008706"
Set-CellText $ws 17 5 "  -11.12%  "
Set-CellText $ws 18 4 "27.922.22"
Set-CellText $ws 18 5 "  -3.74%  "
Set-CellText $ws 19 4 "5.313"
Set-CellText $ws 19 5 "  -4.05%  "
Set-CellText $ws 20 5 "  -0.17%  "
Set-CellText $ws 21 4 "205.36"
Set-CellText $ws 21 5 "  -4.77%  "
Set-CellText $ws 22 5 "  -2.36%  "
Set-CellText $ws 23 4 "6.634"
Set-CellText $ws 23 5 "  -2.95%  "
Set-CellText $ws 24 4 "1.003"
Set-CellText $ws 24 5 "  -0.13%  "
Set-CellText $ws 25 4 "150.28"
Set-CellText $ws 25 5 "  -3.60%  "
Set-CellText $ws 26 4 "8.066"
Set-CellText $ws 26 5 "  +1.44%  "
Set-CellText $ws 27 5 "  -4.12%  "
Set-CellText $ws 29 4 "0.06206"
Set-CellText $ws 29 5 "  -3.47%  "
Set-CellText $ws 30 4 "1.384"
Set-CellText $ws 30 5 "  -3.07%  "
Set-CellText $ws 31 4 "1.394"
Set-CellText $ws 31 5 "  -3.19%  "
Set-CellText $ws 32 4 "3.745"
Set-CellText $ws 32 5 "  -1.69%  "
Set-CellText $ws 33 4 "3.739"
Set-CellText $ws 33 5 "  -0.95%  "
Set-CellText $ws 34 4 "1.680"
Set-CellText $ws 34 5 "  -2.19%  "
Set-CellText $ws 35 4 "1.038"
Set-CellText $ws 35 5 "  -4.82%  "
Set-CellText $ws 36 4 "0.6417"
Set-CellText $ws 36 5 "  -0.41%  "
Set-CellText $ws 37 4 "2.444"
Set-CellText $ws 37 5 "  -3.86%  "
Set-CellText $ws 38 5 "  -1.34%  "
Set-CellText $ws 39 4 "0.01670"
Set-CellText $ws 39 5 "  -4.56%  "
Set-CellText $ws 40 4 "1.122.93"
Set-CellText $ws 40 5 "  -0.70%  "
Set-CellText $ws 41 4 "6.161"
Set-CellText $ws 41 5 "  -6.67%  "
Set-CellText $ws 42 4 "0.8726"
Set-CellText $ws 42 5 "  -1.46%  "
Set-CellText $ws 43 5 "  +0.11%  "
Set-CellText $ws 44 4 "100.09"
Set-CellText $ws 44 5 "  -0.63%  "
Set-CellText $ws 45 4 "1.887.33"
Set-CellText $ws 45 5 "  -4.85%  "
Set-CellText $ws 46 4 "59.40"
Set-CellText $ws 46 5 "  -4.10%  "
Set-CellText $ws 47 5 "  -3.82%  "
Set-CellText $ws 48 4 "1.578"
Set-CellText $ws 48 5 "  -2.06%  "
Set-CellText $ws 49 4 "8.264"
Set-CellText $ws 49 5 "  -1.81%  "
Set-CellText $ws 50 4 "0.05380"
Set-CellText $ws 50 5 "  -2.16%  "
Set-CellText $ws 51 4 "0.4414"
Set-CellText $ws 51 5 "  -2.60%  "
